$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.265
$ws.Range("D3").Value = -7.479000000000001
$ws.Range("A12").Value = -21.519
$ws.Range("B14").Value = 6.233
$ws.Range("D20").Value = -7.765000000000001
$ws.Range("D25").Value = -7.855
$ws.Range("B26").Value = 6.225
$ws.Range("A27").Value = -21.424
$ws.Range("D30").Value = -7.140000000000001
$ws.Range("B31").Value = 6.425
$ws.Range("A32").Value = -21.304
$ws.Range("B35").Value = 7.287000000000001
$ws.Range("A36").Value = -21.112
$ws.Range("B37").Value = 7.848000000000002
$ws.Range("A38").Value = -20.093
$ws.Range("D44").Value = -7.925
$ws.Range("B45").Value = 5.894
$ws.Range("A46").Value = -21.481
$ws.Range("D47").Value = -7.576000000000001
$ws.Range("B52").Value = 5.359
$ws.Range("A54").Value = -21.862
$ws.Range("A55").Value = -22.21
$ws.Range("A56").Value = -22.097
$ws.Range("B57").Value = 5.331999999999999
$ws.Range("D58").Value = -8.178999999999998
$ws.Range("A67").Value = -21.586
$ws.Range("A69").Value = -21.637
$ws.Range("A72").Value = -21.567
$ws.Range("D78").Value = -7.812
$ws.Range("B81").Value = 6.375999999999999
$ws.Range("A83").Value = -20.146
$ws.Range("B83").Value = 7.326000000000001
$ws.Range("D84").Value = -8.196000000000002
$ws.Range("A86").Value = -22.096
$ws.Range("D89").Value = -6.808
$ws.Range("A91").Value = -21.48
$ws.Range("D91").Value = -6.915000000000001
$ws.Range("D92").Value = -6.754
$ws.Range("A93").Value = -21.665
$ws.Range("D96").Value = -7.511
$ws.Range("A99").Value = -20.437
$ws.Range("B100").Value = 5.517
$ws.Range("B102").Value = 7.499000000000001
$ws.Range("D102").Value = -7.865
